$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Game" -> "Games" (shared string already used by B3, edited in place)
$ws.Range("B3").Value = "Games"

# C2 / C3 now hold text values ("1" and "2020.") instead of numbers.
# Use the leading apostrophe so Excel stores them as text rather than
# re-interpreting the numeric-looking text back into a number, then
# clear the resulting "number stored as text" formatting flag so the
# cells keep their original (default) style.
$ws.Range("C2").Value = "'1"
$ws.Range("C3").Value = "'2020."
$ws.Range("C2:C3").ClearFormats()
